# "Fruta / hortaliza, semanal" — weekly refresh: insert the two newest
# price-report rows for this market/variety block (dated 2023-02-07) at the
# top of the block (row 378), pushing the older rows (previously 378-482)
# down by two rows to 380-484. The sheet dimension grows from A1:R482 to
# A1:R484 automatically via the row insert.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 378, shifting the existing data (old rows
# 378-482) down to 380-484.
$ws.Rows("378:379").Insert()

# New row 378 — "Primera" quality, newest report date.
$ws.Range("A378").Value = 1
$ws.Range("B378").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C378").Value = 'Arica y Parinacota'
$ws.Range("D378").Value = 44964
$ws.Range("E378").Value = 15
$ws.Range("F378").Value = 100112032
$ws.Range("G378").Value = 'Zapallo italiano'
$ws.Range("H378").Value = 'Huracán'
$ws.Range("I378").Value = 'Primera'
$ws.Range("J378").Value = 130
$ws.Range("K378").Value = 2500
$ws.Range("L378").Value = 3000
$ws.Range("M378").Value = 2750
$ws.Range("N378").Value = '$/caja 70 unidades'
$ws.Range("O378").Value = 'Región de Arica y Parinacota'
$ws.Range("P378").Value = 39
$ws.Range("Q378").Value = 70
$ws.Range("R378").Value = 'Hortaliza'

# New row 379 — "Segunda" quality, same newest report date.
$ws.Range("A379").Value = 1
$ws.Range("B379").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C379").Value = 'Arica y Parinacota'
$ws.Range("D379").Value = 44964
$ws.Range("E379").Value = 15
$ws.Range("F379").Value = 100112032
$ws.Range("G379").Value = 'Zapallo italiano'
$ws.Range("H379").Value = 'Huracán'
$ws.Range("I379").Value = 'Segunda'
$ws.Range("J379").Value = 130
$ws.Range("K379").Value = 2000
$ws.Range("L379").Value = 2500
$ws.Range("M379").Value = 2250
$ws.Range("N379").Value = '$/caja 100 unidades'
$ws.Range("O379").Value = 'Región de Arica y Parinacota'
$ws.Range("P379").Value = 22
$ws.Range("Q379").Value = 100
$ws.Range("R379").Value = 'Hortaliza'
